# Insert a new data row right above the current row 41 ("Vega Modelo de
# Temuco" / "Bruselas (repollito)" price sheet). Excel shifts every
# following row down by one (old row 41 -> 42, ..., old row 98 -> 99) and
# the sheet's used range grows from A1:R98 to A1:R99 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new price record.
$ws.Cells.Item(41, 1).Value  = 10
$ws.Cells.Item(41, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(41, 3).Value  = "La Araucanía"
$ws.Cells.Item(41, 4).Value  = 44757
$ws.Cells.Item(41, 5).Value  = 9
$ws.Cells.Item(41, 6).Value  = 100112035
$ws.Cells.Item(41, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(41, 8).Value  = "Sin especificar"
$ws.Cells.Item(41, 9).Value  = "Primera"
$ws.Cells.Item(41, 10).Value = 140
$ws.Cells.Item(41, 11).Value = 25000
$ws.Cells.Item(41, 12).Value = 26000
$ws.Cells.Item(41, 13).Value = 25571
$ws.Cells.Item(41, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(41, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(41, 16).Value = 2557
$ws.Cells.Item(41, 17).Value = 10
$ws.Cells.Item(41, 18).Value = "Hortaliza"
